$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Locoto" (Vega Modelo de Temuco).
# It is inserted as a new row 91, pushing the existing rows 91-102 down to
# 92-103 (matches the diff: row 91 becomes new data, old rows 91-101 shift
# to 92-102, and old row 102 becomes the new row 103).
$ws.Rows.Item(91).Insert()

$ws.Range("A91").Value = 10
$ws.Range("B91").Value = 'Vega Modelo de Temuco'
$ws.Range("C91").Value = 'La Araucanía'
$ws.Range("D91").Value = 45180
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = 100112042
$ws.Range("G91").Value = 'Locoto'
$ws.Range("H91").Value = 'Sin especificar'
$ws.Range("I91").Value = 'Primera'
$ws.Range("J91").Value = 80
$ws.Range("K91").Value = 2200
$ws.Range("L91").Value = 2200
$ws.Range("M91").Value = 2200
$ws.Range("N91").Value = '$/kilo'
$ws.Range("O91").Value = 'Región de Arica y Parinacota'
$ws.Range("P91").Value = 2200
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = 'Hortaliza'
